# Auto-generated Excel COM-interop script
# Applies updated market-price data (columns H-N) to specific leve rows
# across multiple sheets, per the commit's scheduled data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 306599.88
$ws.Range("J17").Value = 314931.12
$ws.Range("L17").Value = 944793.36
$ws.Range("N17").Value = -945129.36
$ws.Range("H18").Value = 609.5625
$ws.Range("I18").Value = 614.74194
$ws.Range("K18").Value = 614.74194
$ws.Range("M18").Value = -330.74194
$ws.Range("H88").Value = 1545.8889
$ws.Range("J88").Value = 1733.75
$ws.Range("L88").Value = 1733.75
$ws.Range("N88").Value = -2545.75
$ws.Range("H91").Value = 1545.8889
$ws.Range("J91").Value = 1733.75
$ws.Range("L91").Value = 1733.75
$ws.Range("N91").Value = -4541.75
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H97").Value = 24625024
$ws.Range("J97").Value = 30781154
$ws.Range("L97").Value = 92343462
$ws.Range("N97").Value = -92344454
$ws.Range("H101").Value = 935.26666
$ws.Range("I101").Value = 439.8
$ws.Range("K101").Value = 1319.4
$ws.Range("M101").Value = 302.5999999999999
$ws.Range("H137").Value = 21277526
$ws.Range("I137").Value = 30304024
$ws.Range("J137").Value = 779.3570999999999
$ws.Range("K137").Value = 90912072
$ws.Range("L137").Value = 2338.0713
$ws.Range("M137").Value = -90909522
$ws.Range("N137").Value = -7438.0713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 373548.34
$ws.Range("I45").Value = 795573.1
$ws.Range("J45").Value = 4276.625
$ws.Range("K45").Value = 795573.1
$ws.Range("L45").Value = 4276.625
$ws.Range("M45").Value = -795196.1
$ws.Range("N45").Value = -5030.625
$ws.Range("H97").Value = 551.13794
$ws.Range("I97").Value = 570.88464
$ws.Range("J97").Value = 380
$ws.Range("K97").Value = 570.88464
$ws.Range("L97").Value = 380
$ws.Range("M97").Value = -74.88463999999999
$ws.Range("N97").Value = -1372
$ws.Range("H132").Value = 892.0952
$ws.Range("I132").Value = 848.1842
$ws.Range("K132").Value = 2544.5526
$ws.Range("M132").Value = -14.55259999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 49998.5
$ws.Range("J76").Value = 49998.5
$ws.Range("L76").Value = 49998.5
$ws.Range("N76").Value = -50628.5
$ws.Range("H79").Value = 49998.5
$ws.Range("J79").Value = 49998.5
$ws.Range("L79").Value = 49998.5
$ws.Range("N79").Value = -52182.5
$ws.Range("H107").Value = 100001500
$ws.Range("I107").Value = 1510
$ws.Range("J107").Value = 166668180
$ws.Range("K107").Value = 1510
$ws.Range("L107").Value = 166668180
$ws.Range("M107").Value = 410
$ws.Range("N107").Value = -166672020

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4488.5
$ws.Range("I31").Value = 5557.448
$ws.Range("K31").Value = 5557.448
$ws.Range("M31").Value = -5262.448
$ws.Range("H34").Value = 4488.5
$ws.Range("I34").Value = 5557.448
$ws.Range("K34").Value = 5557.448
$ws.Range("M34").Value = -5355.448
$ws.Range("H122").Value = 2363.6072
$ws.Range("I122").Value = 2591.842
$ws.Range("J122").Value = 1881.7778
$ws.Range("K122").Value = 7775.526
$ws.Range("L122").Value = 5645.3334
$ws.Range("M122").Value = -5325.526
$ws.Range("N122").Value = -10545.3334
$ws.Range("H132").Value = 2812.1714
$ws.Range("I132").Value = 2727.862
$ws.Range("K132").Value = 8183.586
$ws.Range("M132").Value = -5653.586
$ws.Range("H134").Value = 1463.2916
$ws.Range("I134").Value = 1278.2727
$ws.Range("J134").Value = 3498.5
$ws.Range("K134").Value = 3834.8181
$ws.Range("L134").Value = 10495.5
$ws.Range("M134").Value = -1299.8181
$ws.Range("N134").Value = -15565.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 699.8889
$ws.Range("I13").Value = 322.25
$ws.Range("J13").Value = 1002
$ws.Range("K13").Value = 966.75
$ws.Range("L13").Value = 3006
$ws.Range("M13").Value = -798.75
$ws.Range("N13").Value = -3342
$ws.Range("H132").Value = 1826
$ws.Range("I132").Value = 1995.6666
$ws.Range("J132").Value = 1698.75
$ws.Range("K132").Value = 17960.9994
$ws.Range("L132").Value = 15288.75
$ws.Range("M132").Value = -15430.9994
$ws.Range("N132").Value = -20348.75
$ws.Range("H134").Value = 1816.7894
$ws.Range("I134").Value = 1362.1666
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 4086.4998
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = 983.5001999999999
$ws.Range("N134").Value = -40140
$ws.Range("H137").Value = 1878.9375
$ws.Range("I137").Value = 1306.4
$ws.Range("J137").Value = 2833.1667
$ws.Range("K137").Value = 3919.2
$ws.Range("L137").Value = 8499.500100000001
$ws.Range("M137").Value = 1180.8
$ws.Range("N137").Value = -18699.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 1500
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 2000
$ws.Range("K29").Value = 1000
$ws.Range("L29").Value = 2000
$ws.Range("M29").Value = -710
$ws.Range("N29").Value = -2580
$ws.Range("H43").Value = 3120.9333
$ws.Range("J43").Value = 5799.8
$ws.Range("L43").Value = 5799.8
$ws.Range("N43").Value = -6101.8
$ws.Range("H95").Value = 34842.57
$ws.Range("J95").Value = 35316.332
$ws.Range("L95").Value = 35316.332
$ws.Range("N95").Value = -40808.332
$ws.Range("H102").Value = 1835
$ws.Range("I102").Value = 1470.2
$ws.Range("K102").Value = 1470.2
$ws.Range("M102").Value = 151.8
$ws.Range("H132").Value = 3023.46
$ws.Range("I132").Value = 2544
$ws.Range("J132").Value = 4541.75
$ws.Range("K132").Value = 7632
$ws.Range("L132").Value = 13625.25
$ws.Range("M132").Value = -5102
$ws.Range("N132").Value = -18685.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2248.5
$ws.Range("I40").Value = 1729.6818
$ws.Range("J40").Value = 4150.8335
$ws.Range("K40").Value = 1729.6818
$ws.Range("L40").Value = 4150.8335
$ws.Range("M40").Value = -1593.6818
$ws.Range("N40").Value = -4422.8335
$ws.Range("H55").Value = 561.95654
$ws.Range("I55").Value = 402.26666
$ws.Range("J55").Value = 861.375
$ws.Range("K55").Value = 402.26666
$ws.Range("L55").Value = 861.375
$ws.Range("M55").Value = -229.26666
$ws.Range("N55").Value = -1207.375
$ws.Range("H132").Value = 4750
$ws.Range("J132").Value = 11671.75
$ws.Range("L132").Value = 35015.25
$ws.Range("N132").Value = -40075.25
$ws.Range("H136").Value = 2224.3057
$ws.Range("I136").Value = 1877.75
$ws.Range("K136").Value = 5633.25
$ws.Range("M136").Value = -3083.25
$ws.Range("H140").Value = 24166.334
$ws.Range("J140").Value = 15999.75
$ws.Range("L140").Value = 15999.75
$ws.Range("N140").Value = -26359.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1708064.9
$ws.Range("J62").Value = 8310.75
$ws.Range("L62").Value = 8310.75
$ws.Range("N62").Value = -9558.75
$ws.Range("H65").Value = 1708064.9
$ws.Range("J65").Value = 8310.75
$ws.Range("L65").Value = 41553.75
$ws.Range("N65").Value = -47793.75
$ws.Range("H70").Value = 33023.75
$ws.Range("J70").Value = 35666.668
$ws.Range("L70").Value = 35666.668
$ws.Range("N70").Value = -36296.668
$ws.Range("H73").Value = 33023.75
$ws.Range("J73").Value = 35666.668
$ws.Range("L73").Value = 35666.668
$ws.Range("N73").Value = -37850.668
$ws.Range("H126").Value = 1705.9286
$ws.Range("I126").Value = 1557.3334
$ws.Range("K126").Value = 4672.0002
$ws.Range("M126").Value = -2202.0002
$ws.Range("H132").Value = 1825.3448
$ws.Range("I132").Value = 1639.62
$ws.Range("K132").Value = 4918.86
$ws.Range("M132").Value = -2388.86
